# "upload of rainfall and yield data"
#
# The sheet had a two-row header (row 2 = column titles, row 3 = unit
# sub-titles "(Rai)" / "(Rai)" / "(tons)") with A2:A3 merged for the
# "Year1" label. This edit removes the unit sub-header row entirely -
# the data rows (1981-2016) shift up by one and the now-unused "(Rai)"
# and "(tons)" shared-string entries drop out of the saved workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the sub-header row (old row 3: "(Rai)", "(Rai)", "(tons)").
# This also removes the A2:A3 merged cell and shifts every row below
# it (the 1981-2016 data, previously rows 4-39) up by one, so the
# dimension becomes A1:F38 instead of A1:F39.
$ws.Rows(3).Delete()
